$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.07822233795641864
$ws.Range("B1").Value = 0.0778929003770088
$ws.Range("A2").Value = 0.031272088310627666
$ws.Range("B2").Value = -0.0319210655582598
$ws.Range("A3").Value = 0.16034575632048842
$ws.Range("B3").Value = -0.16107651537648238
$ws.Range("A4").Value = -0.12291883288645167
$ws.Range("B4").Value = 0.1224435237023549
$ws.Range("A5").Value = -0.11644352391616941
$ws.Range("B5").Value = 0.11551494290603159
$ws.Range("A6").Value = -0.07101759473760971
$ws.Range("B6").Value = 0.07093707876593891
$ws.Range("A7").Value = -0.05093707902707045
$ws.Range("B7").Value = 0.05077561378645257
$ws.Range("A8").Value = -0.030775614049552757
$ws.Range("B8").Value = 0.030684066090749695
$ws.Range("A9").Value = -0.024684066314523356
$ws.Range("B9").Value = 0.024619167735600023
$ws.Range("A10").Value = -0.018619167960864047
$ws.Range("B10").Value = 0.018611950075914763
$ws.Range("A11").Value = -0.01411195029643153
$ws.Range("B11").Value = 0.01410369516433363
$ws.Range("A12").Value = -0.01999745388366847
$ws.Range("B12").Value = 0.019954735895569176
$ws.Range("A13").Value = -0.013954736122401279
$ws.Range("B13").Value = 0.013944273701665288
$ws.Range("A14").Value = -0.0019442739465329595
$ws.Range("B14").Value = 0.001931553646908668
$ws.Range("A15").Value = 0.004068446125788938
$ws.Range("B15").Value = -0.004087114139217363
$ws.Range("A16").Value = -0.015026565410077186
$ws.Range("B16").Value = 0.015004388730808227
$ws.Range("A17").Value = -0.00900438895900546
$ws.Range("B17").Value = 0.008999999762704824
$ws.Range("A18").Value = -0.09916498891420034
$ws.Range("B18").Value = 0.09900007035671976
$ws.Range("A19").Value = -0.09000007058116566
$ws.Range("B19").Value = 0.08869110413320458
$ws.Range("A20").Value = -0.01801427736922001
$ws.Range("B20").Value = 0.01800437646419706
$ws.Range("A21").Value = -0.00900437669639409
$ws.Range("B21").Value = 0.008999999767531186
$ws.Range("A22").Value = -0.11942469098531383
$ws.Range("B22").Value = 0.11892420860936781
$ws.Range("A23").Value = -0.10992420883516463
$ws.Range("B23").Value = 0.10907573633796286
$ws.Range("A24").Value = -0.042125447349862455
$ws.Range("B24").Value = 0.0419999996700362
$ws.Range("A25").Value = -0.03849734838766494
$ws.Range("B25").Value = 0.038473154992868075
$ws.Range("A26").Value = -0.032473155212390026
$ws.Range("B26").Value = 0.032447639344251655
$ws.Range("A27").Value = -0.02644763956419638
$ws.Range("B27").Value = 0.026377846841088903
$ws.Range("A28").Value = -0.020377847062522214
$ws.Range("B28").Value = 0.02034307738720198
$ws.Range("A29").Value = -0.008343077627118944
$ws.Range("B29").Value = 0.008338792145895724
$ws.Range("A30").Value = -0.04215944649679981
$ws.Range("B30").Value = 0.04201903449803979
$ws.Range("A31").Value = -0.027019034749162785
$ws.Range("B31").Value = 0.027000854655083728
$ws.Range("A32").Value = -0.006000854924282706
$ws.Range("B32").Value = 0.005999999774831011
